# Fruta / hortaliza, semanal
# The rows 2-9 (date D, volume M, min/max/avg price N/O/P and $/Kg price S)
# get shuffled among themselves (a permutation of the weekly records),
# while all other columns (A,B,C,E,F,G,H,I,J,K,L,Q,R,T) stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get reshuffled between rows.
$cols = @("D", "M", "N", "O", "P", "S")

# Snapshot current values for the affected columns, rows 2-9.
# Use Value2() to get the raw underlying number (avoids Variant/DateTime
# conversion quirks tied to the date number format on column D).
$orig = @{}
foreach ($r in 2..9) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2()
    }
    $orig[$r] = $rowVals
}

# New row r gets the values that used to live in row $map[r].
$map = @{
    2 = 3
    3 = 5
    4 = 6
    5 = 9
    6 = 4
    7 = 8
    8 = 7
    9 = 2
}

foreach ($r in 2..9) {
    $src = $map[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $orig[$src][$c]
    }
}
